# Notes.xlsx conversion fix-up:
#  - swap the FirstName/LastName header labels in B1/C1
#  - renumber the CNE id column (A2:A11) to the corrected values
#  - re-apply formatting (font color) across A1:C11
#  - restore the view scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: B1/C1 swap (FirstName <-> LastName) ---------------------
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# --- Column A (CNE) renumbering for the data rows -------------------------
$ids = @(18000021, 18000022, 18000023, 18000024, 18000025, 18000026, 18000027, 18000028, 18000029, 18000030)
for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# --- Formatting: re-apply font color across the whole A1:C11 block -------
$ws.Range("A1:C11").Font.ThemeColor = 1

# --- View state: scroll back to column A, select G10 ----------------------
[void]$ws.Activate()
[void]$ws.Range("A1").Select()
[void]$ws.Range("G10").Select()
